$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update student record in row 2 with real name info, replacing the
# sample placeholder data (unactive role / user account update).
# New shared-string entries must land in this order: Benito, Cortes, Masubong.
$ws.Range("C2").Value = "Benito"
$ws.Range("D2").Value = "Cortes"
$ws.Range("B2").Value = "Masubong"

# Move the active selection to J3 (matches the saved cursor position).
$ws.Range("J3").Select()
